$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.565.74'
$ws.Cells.Item(2, 5).Value = '  +4.33%  '
$ws.Cells.Item(3, 4).Value = '1.792.29'
$ws.Cells.Item(3, 5).Value = '  +0.70%  '
$ws.Cells.Item(4, 5).Value = '  +0.16%  '
$ws.Cells.Item(5, 4).Value = '''313.79'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.04%  '
$ws.Cells.Item(7, 4).Value = '''0.5360'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.17%  '
$ws.Cells.Item(8, 4).Value = '''0.3807'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.71%  '
$ws.Cells.Item(9, 4).Value = '''0.07521'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +1.64%  '
$ws.Cells.Item(10, 4).Value = '''42.47'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -1.04%  '
$ws.Cells.Item(11, 4).Value = '''1.116'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.83%  '
$ws.Cells.Item(12, 4).Value = '''1.002'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.20%  '
$ws.Cells.Item(13, 4).Value = '''21.10'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +2.16%  '
$ws.Cells.Item(14, 4).Value = '''6.181'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.24%  '
$ws.Cells.Item(15, 4).Value = '''7.424'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +6.24%  '
$ws.Cells.Item(16, 4).Value = '1.790.42'
$ws.Cells.Item(17, 4).Value = '''90.21'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.73%  '
$ws.Cells.Item(18, 4).Value = '''0.00001065'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.88%  '
$ws.Cells.Item(19, 4).Value = '''0.06444'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(20, 4).Value = '''1.001'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.12%  '
$ws.Cells.Item(21, 4).Value = '''17.24'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +2.57%  '
$ws.Cells.Item(22, 4).Value = '''5.921'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.31%  '
$ws.Cells.Item(23, 4).Value = '28.574.86'
$ws.Cells.Item(23, 5).Value = '  +4.20%  '
$ws.Cells.Item(24, 4).Value = '''11.21'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.22%  '
$ws.Cells.Item(25, 4).Value = '''2.098'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.79%  '
$ws.Cells.Item(26, 4).Value = '''161.25'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +3.54%  '
$ws.Cells.Item(27, 4).Value = '''20.50'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +1.70%  '
$ws.Cells.Item(28, 4).Value = '''2.381'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.92%  '
$ws.Cells.Item(29, 4).Value = '1.996.64'
$ws.Cells.Item(29, 5).Value = '  +0.56%  '
$ws.Cells.Item(30, 4).Value = '''123.16'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +1.79%  '
$ws.Cells.Item(31, 4).Value = '''1.122'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +3.92%  '
$ws.Cells.Item(32, 4).Value = '''0.1021'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.99%  '
$ws.Cells.Item(33, 4).Value = '''5.707'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +2.03%  '
$ws.Cells.Item(34, 4).Value = '''3.657'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.98%  '
$ws.Cells.Item(35, 4).Value = '''0.2306'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +11.94%  '
$ws.Cells.Item(36, 4).Value = '''0.06552'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +9.65%  '
$ws.Cells.Item(37, 4).Value = '''0.02322'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +2.78%  '
$ws.Cells.Item(38, 2).Value = 'FraxShare'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(38, 4).Value = '''8.706'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +5.21%  '
$ws.Cells.Item(39, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(39, 4).Value = '''5.088'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +3.36%  '
$ws.Cells.Item(40, 4).Value = '''11.46'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +1.83%  '
$ws.Cells.Item(41, 4).Value = '''0.6317'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +3.38%  '
$ws.Cells.Item(42, 4).Value = '''1.213'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +6.94%  '
$ws.Cells.Item(44, 4).Value = '''1.382'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -3.07%  '
$ws.Cells.Item(45, 4).Value = '''13.57'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +1.84%  '
$ws.Cells.Item(46, 4).Value = '''0.5921'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +2.59%  '
$ws.Cells.Item(47, 4).Value = '''3.667'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +1.14%  '
$ws.Cells.Item(48, 4).Value = '''125.28'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +3.42%  '
$ws.Cells.Item(49, 4).Value = '''1.979'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +4.69%  '
$ws.Cells.Item(50, 5).Value = '  +4.16%  '
$ws.Cells.Item(51, 4).Value = '''0.06928'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +3.00%  '
